$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Saudi Arabia" rows (one per "along" group: Global, High-income,
# International). Deleting bottom-up keeps the earlier row numbers valid and
# shifts the remaining rows (and their shared-string E refs) up automatically.
$ws.Rows.Item(39).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(13).Delete()

# Re-run values for RU 1001 ("without crop") changed the aggregated "All" and
# "Europe" rows (which include Russia) plus the "Russia" row itself, in each
# of the three remaining groups. Every other country row is unaffected.

# Global group
$ws.Range("B2").Value = 73.8161485502559
$ws.Range("C2").Value = 72.4542861110033
$ws.Range("D2").Value = 75.1780109895084

$ws.Range("C3").Value = 75.6462846457027
$ws.Range("D3").Value = 79.8462128013628

$ws.Range("B12").Value = 76.7816309654674
$ws.Range("C12").Value = 72.3067698657709
$ws.Range("D12").Value = 81.2564920651638

# High-income group
$ws.Range("B14").Value = 69.211130206042
$ws.Range("C14").Value = 67.7830711918507
$ws.Range("D14").Value = 70.6391892202333

$ws.Range("C15").Value = 69.1652730389077
$ws.Range("D15").Value = 73.7345785221246

$ws.Range("B24").Value = 69.4644695949361
$ws.Range("C24").Value = 64.5975635112594
$ws.Range("D24").Value = 74.3313756786127

# International group
$ws.Range("B26").Value = 68.2106860437125
$ws.Range("C26").Value = 66.767714507923
$ws.Range("D26").Value = 69.6536575795019

$ws.Range("C27").Value = 69.2393895182127
$ws.Range("D27").Value = 73.8436649948966

$ws.Range("B36").Value = 73.799243065225
$ws.Range("C36").Value = 68.9277869442202
$ws.Range("D36").Value = 78.6706991862297
